# Apply the edit described by the diff:
# 1. Rename the 9 worksheets from their old "summ*" names to the new
#    randomized-suffix "summ*" names (same order as they appear in the workbook).
# 2. On every sheet, change cell A31 from "CarAvailable" to "CarOwnershipHH".

$wb = $excel.ActiveWorkbook

# Mapping of old sheet name -> new sheet name, in sheet (tab) order.
$newNames = @(
    "summ22990037",
    "summ16653750",
    "summ10588009",
    "summ04739417",
    "summ02386118",
    "summ59957307",
    "summ00290844",
    "summ58177010",
    "summ59384591"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Update the label in column A, row 31.
    if ($ws.Range("A31").Value() -eq "CarAvailable") {
        $ws.Range("A31").Value = "CarOwnershipHH"
    }

    # Rename the worksheet tab.
    $ws.Name = $newNames[$i - 1]
}
